# DataPengiriman.xlsx -- "fixing data pengiriman excel dengan validasi
# format no hp" (fix phone-number formatting in the shipment data).
#
# The recipient/sender phone numbers in columns I (no_hp_pengirim) and J
# (no_hp_penerima) were stored without the Indonesian country code, so they
# get normalised to the full "62" + number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- phone number corrections (no_hp_pengirim = I, no_hp_penerima = J) ---
$ws.Range("I2").Value = 6281627364523
$ws.Range("J2").Value = 6281726354728

$ws.Range("I3").Value = 6287564637261
$ws.Range("J3").Value = 6281274657487

$ws.Range("I4").Value = 6281254647587
$ws.Range("J4").Value = 6287364758676

$ws.Range("I5").Value = 628717267643
$ws.Range("J5").Value = 6282736457485

$ws.Range("I6").Value = 6281672537485
$ws.Range("J6").Value = 628172645362

# --- widen column I a touch so the longer numbers remain fully visible ---
$ws.Range("I1").ColumnWidth = 15.6

# --- view state left behind by Excel on save (zoom + active selection) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 75
$ws.Range("J7").Select()
